# Weekly price update: insert a new record for the week of 2022-01-24
# (serial 44585) as row 175, shifting the existing rows 175-184 down to
# 176-185 (matching the canonical diff where the whole block of rows from
# 175 onward is pushed down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 175, pushing rows 175..184 to 176..185.
$ws.Rows.Item(175).Insert()

# Populate the newly inserted row 175 with the new weekly record.
$ws.Cells.Item(175, 1).Value = 4
$ws.Cells.Item(175, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(175, 3).Value = "Los Lagos"
$ws.Cells.Item(175, 4).Value = 44585
$ws.Cells.Item(175, 5).Value = 10
$ws.Cells.Item(175, 6).Value = 100112032
$ws.Cells.Item(175, 7).Value = "Zapallo italiano"
$ws.Cells.Item(175, 8).Value = "Sin especificar"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 80
$ws.Cells.Item(175, 11).Value = 16000
$ws.Cells.Item(175, 12).Value = 16000
$ws.Cells.Item(175, 13).Value = 16000
$ws.Cells.Item(175, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(175, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(175, 16).Value = 320
$ws.Cells.Item(175, 17).Value = 50
$ws.Cells.Item(175, 18).Value = "Hortaliza"
